$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("D2").Value = 28
$ws.Range("G2").Value = 40
$ws.Range("H2").Value = 41
$ws.Range("E3").Value = 50
$ws.Range("H3").Value = 35
$ws.Range("I3").Value = 71
$ws.Range("C9").Value = 173
$ws.Range("D9").Value = 161
$ws.Range("E9").Value = 167
$ws.Range("G9").Value = 194
$ws.Range("I9").Value = 204
$ws.Range("B10").Value = 392
$ws.Range("C10").Value = 472
$ws.Range("D10").Value = 628
$ws.Range("E10").Value = 705
$ws.Range("F10").Value = 829
$ws.Range("G10").Value = 498
$ws.Range("H10").Value = 164
$ws.Range("I10").Value = 297
$ws.Range("J10").Value = 260
$ws.Range("B11").Value = 581
$ws.Range("C11").Value = 704
$ws.Range("D11").Value = 874
$ws.Range("E11").Value = 953
$ws.Range("F11").Value = 1098
$ws.Range("G11").Value = 780
$ws.Range("H11").Value = 393
$ws.Range("I11").Value = 615
$ws.Range("J11").Value = 537

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("E7").Value = 26
$ws.Range("E8").Value = 42

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("D2").Value = 3
$ws.Range("G2").Value = 3
$ws.Range("F7").Value = 13
$ws.Range("D8").Value = 27
$ws.Range("F8").Value = 26
$ws.Range("G8").Value = 16

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("H3").Value = 6
$ws.Range("I7").Value = 35
$ws.Range("B8").Value = 48
$ws.Range("D8").Value = 188
$ws.Range("E8").Value = 206
$ws.Range("F8").Value = 252
$ws.Range("H8").Value = 17
$ws.Range("I8").Value = 78
$ws.Range("B9").Value = 63
$ws.Range("D9").Value = 222
$ws.Range("E9").Value = 244
$ws.Range("F9").Value = 280
$ws.Range("H9").Value = 47
$ws.Range("I9").Value = 129

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 7

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("D6").Value = 16
$ws.Range("F6").Value = 45
$ws.Range("D7").Value = 20
$ws.Range("F7").Value = 51

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I5").Value = 4
$ws.Range("J6").Value = 8
$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I6").Value = 9
$ws.Range("I8").Value = 15

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 3
$ws.Range("C5").Value = 7
$ws.Range("B7").Value = 17
$ws.Range("G8").Value = 61
$ws.Range("G14").Value = 5
$ws.Range("C28").Value = 48
$ws.Range("F28").Value = 35
$ws.Range("I28").Value = 28
$ws.Range("E29").Value = 7
$ws.Range("C30").Value = 7
$ws.Range("D30").Value = 18
$ws.Range("E32").Value = 42
$ws.Range("D36").Value = 27
$ws.Range("F36").Value = 26
$ws.Range("G36").Value = 16
$ws.Range("F38").Value = 7
$ws.Range("C41").Value = 9
$ws.Range("B45").Value = 7
$ws.Range("C47").Value = 23
$ws.Range("E47").Value = 17
$ws.Range("H47").Value = 15
$ws.Range("I47").Value = 11
$ws.Range("D49").Value = 4
$ws.Range("I50").Value = 6
$ws.Range("J50").Value = 15
$ws.Range("C52").Value = 11
$ws.Range("B53").Value = 63
$ws.Range("D53").Value = 222
$ws.Range("E53").Value = 244
$ws.Range("F53").Value = 280
$ws.Range("H53").Value = 47
$ws.Range("I53").Value = 129
$ws.Range("C61").Value = 12
$ws.Range("I65").Value = 15
$ws.Range("G66").Value = 2
$ws.Range("D70").Value = 20
$ws.Range("F70").Value = 51
$ws.Range("D74").Value = 26
$ws.Range("D76").Value = 16
$ws.Range("C78").Value = 7
$ws.Range("E80").Value = 7
$ws.Range("H86").Value = 11
$ws.Range("E94").Value = 35
$ws.Range("I97").Value = 2
$ws.Range("B98").Value = 581
$ws.Range("C98").Value = 704
$ws.Range("D98").Value = 874
$ws.Range("E98").Value = 953
$ws.Range("F98").Value = 1098
$ws.Range("G98").Value = 780
$ws.Range("H98").Value = 393
$ws.Range("I98").Value = 615
$ws.Range("J98").Value = 537

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("C4").Value = 4
$ws.Range("C6").Value = 9

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("E3").Value = 3
$ws.Range("E6").Value = 7

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("H2").Value = 5
$ws.Range("H7").Value = 11

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 7

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 6
$ws.Range("C5").Value = 33
$ws.Range("F5").Value = 17
$ws.Range("C6").Value = 48
$ws.Range("F6").Value = 35
$ws.Range("I6").Value = 28

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("H3").Value = 2
$ws.Range("C7").Value = 14
$ws.Range("E7").Value = 14
$ws.Range("I7").Value = 5
$ws.Range("C8").Value = 23
$ws.Range("E8").Value = 17
$ws.Range("H8").Value = 15
$ws.Range("I8").Value = 11

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("B5").Value = 7
$ws.Range("B6").Value = 7

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("E5").Value = 1
$ws.Range("E7").Value = 7

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("D7").Value = 13
$ws.Range("D8").Value = 16

$ws = $wb.Worksheets.Item('River North')
$ws.Range("D5").Value = 21
$ws.Range("D6").Value = 26

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("C6").Value = 7
$ws.Range("C7").Value = 11

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("H5").Value = 1
$ws.Range("H7").Value = 3

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("G3").Value = 2
$ws.Range("G5").Value = 5

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("D4").Value = 1
$ws.Range("D6").Value = 4

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("E6").Value = 33
$ws.Range("E7").Value = 35

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("B6").Value = 10
$ws.Range("B7").Value = 17

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("I3").Value = 1
$ws.Range("I7").Value = 2

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 15
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 18

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("F5").Value = 4
$ws.Range("F6").Value = 7

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("G2").Value = 1
$ws.Range("G7").Value = 44
$ws.Range("G8").Value = 61

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 2
